$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# --- Fix E81 / E82: these were stored as text, convert them to real numbers ---
$ws.Range("E81").Value = 533274
$ws.Range("E82").Value = 526371

# --- Append new rows 83-85 with the latest screener results ---

# Row 83 - BSE itself
$ws.Cells.Item(83, 1).Value = "27/06/2024 08:45:10"
$ws.Cells.Item(83, 2).Value = 1
$ws.Cells.Item(83, 3).Value = "BSE"
$ws.Cells.Item(83, 4).Value = "BSE (Bombay stock exchange)"
$ws.Cells.Item(83, 5).NumberFormat = "@"
$ws.Cells.Item(83, 5).Value = "20"
$ws.Cells.Item(83, 6).Value = -0.4
$ws.Cells.Item(83, 7).Value = 2531.05
$ws.Cells.Item(83, 8).Value = 293296

# Row 84 - Prestige Estates Projects Limited
$ws.Cells.Item(84, 1).Value = "27/06/2024 08:45:10"
$ws.Cells.Item(84, 2).Value = 2
$ws.Cells.Item(84, 3).Value = "PRESTIGE"
$ws.Cells.Item(84, 4).Value = "Prestige Estates Projects Limited"
$ws.Cells.Item(84, 5).NumberFormat = "@"
$ws.Cells.Item(84, 5).Value = "533274"
$ws.Cells.Item(84, 6).Value = -3.58
$ws.Cells.Item(84, 7).Value = 1858.15
$ws.Cells.Item(84, 8).Value = 790534

# Row 85 - Nmdc Limited
$ws.Cells.Item(85, 1).Value = "27/06/2024 08:45:10"
$ws.Cells.Item(85, 2).Value = 3
$ws.Cells.Item(85, 3).Value = "NMDC"
$ws.Cells.Item(85, 4).Value = "Nmdc Limited"
$ws.Cells.Item(85, 5).NumberFormat = "@"
$ws.Cells.Item(85, 5).Value = "526371"
$ws.Cells.Item(85, 6).Value = -3.05
$ws.Cells.Item(85, 7).Value = 241.65
$ws.Cells.Item(85, 8).Value = 11962458
